# Add Q4-2022 holding data:
#  - "总计" (summary) sheet gets a new row for 2022-Q4, its old 2022-Q1
#    row moves down to row 3, and a brand-new "2022-Q4" worksheet (holding
#    the quarter's fund detail) is inserted between "总计" and "2022-Q1".

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)   # "总计"

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after the summary sheet.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Header row (B1:H1) - reuse the summary sheet's header formatting.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# A2 - reuse the summary sheet's row-label formatting.
$summary.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)      # xlPasteFormats
$q4.Range("A2").Value = 0

# B2:G2 - plain text fund details (no special styling).
$q4text = $q4.Range("B2:G2")
$q4text.NumberFormat = "@"
$q4.Range("B2").Value = "165531"
$q4.Range("C2").Value = "信诚多策略灵活配置混合（LOF）"
$q4.Range("D2").Value = "0.89"
$q4.Range("E2").Value = "72.25"
$q4.Range("F2").Value = "1.10"
$q4.Range("G2").Value = "0.0098"
$q4text.Style = "Normal"

# H2 - plain numeric rank.
$q4.Range("H2").Value = 3

# Page margins: 0.75/0.75/1/1/0.5/0.5 inch (matches the summary sheet).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: relabel the existing row as
#    2022-Q4 with its real count, and append a new row for 2022-Q1
#    (the data that used to live in row 2).
# ---------------------------------------------------------------------
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.02

# Keep "2022-Q1" as the selected/active tab (it was active before the
# edit, and inserting the new sheet must not steal that selection).
# Looked up by name - worksheet index 2 now refers to the new "2022-Q4"
# sheet since "2022-Q1" shifted down to index 3.
$wb.Worksheets.Item("2022-Q1").Activate()
